$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was logged for this market/product combo. Insert a new
# row at row 119 (pushing the existing rows 119-147 down to 120-148) and
# populate it with the new observation.
$ws.Rows.Item(119).Insert()

$ws.Range("A119").Value = 11
$ws.Range("B119").Value = "Vega Monumental Concepción"
$ws.Range("C119").Value = "Bíobío"
$ws.Range("D119").Value = 44782
$ws.Range("D119").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E119").Value = 8
$ws.Range("F119").Value = 100112043
$ws.Range("G119").Value = "Pepino ensalada"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 100
$ws.Range("K119").Value = 22000
$ws.Range("L119").Value = 23000
$ws.Range("M119").Value = 22500
$ws.Range("N119").Value = "$/caja 60 unidades"
$ws.Range("O119").Value = "Región de Arica y Parinacota"
$ws.Range("P119").Value = 375
$ws.Range("Q119").Value = 60
$ws.Range("R119").Value = "Hortaliza"
